# Weekly refresh of the "Vega Modelo de Temuco - Haba" price series.
# A new weekly price record is inserted as row 63 (pushing the existing
# rows 63-79 down to 64-80), and the worksheet's used-range dimension
# grows from A1:R79 to A1:R80 accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 63, shifting rows 63:79 down to 64:80.
$ws.Rows.Item(63).Insert()

# Populate the new row 63 with the new weekly record.
$ws.Cells.Item(63, 1).Value2  = 10
$ws.Cells.Item(63, 2).Value2  = "Vega Modelo de Temuco"
$ws.Cells.Item(63, 3).Value2  = "La Araucanía"
$ws.Cells.Item(63, 4).Value2  = 44855
$ws.Cells.Item(63, 5).Value2  = 9
$ws.Cells.Item(63, 6).Value2  = 100112026
$ws.Cells.Item(63, 7).Value2  = "Haba"
$ws.Cells.Item(63, 8).Value2  = "Sin especificar"
$ws.Cells.Item(63, 9).Value2  = "Primera"
$ws.Cells.Item(63, 10).Value2 = 70
$ws.Cells.Item(63, 11).Value2 = 9000
$ws.Cells.Item(63, 12).Value2 = 10000
$ws.Cells.Item(63, 13).Value2 = 9571
$ws.Cells.Item(63, 14).Value2 = "`$/saco 25 kilos"
$ws.Cells.Item(63, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(63, 16).Value2 = 383
$ws.Cells.Item(63, 17).Value2 = 25
$ws.Cells.Item(63, 18).Value2 = "Hortaliza"

# Keep the date column formatted the same as the rest of the series.
$ws.Cells.Item(63, 4).NumberFormat = $ws.Cells.Item(64, 4).NumberFormat
